$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list values refreshed (prices + 1h volume %) per GitHub Actions run.
# D-column cells whose new value is a plain decimal number must be forced to
# Text format first, otherwise Excel auto-converts them to a Number and the
# exact textual representation (e.g. trailing zero in "142.40") is lost.

# Row 2
$ws.Cells.Item(2, 4).Value = '59.016.23'
$ws.Cells.Item(2, 5).Value = '  +0.77%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.573.61'
$ws.Cells.Item(3, 5).Value = '  -0.35%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.03%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '562.96'
$ws.Cells.Item(5, 5).Value = '  +3.82%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '142.40'
$ws.Cells.Item(6, 5).Value = '  -1.22%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.02%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +1.84%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.577.91'
$ws.Cells.Item(9, 5).Value = '  -0.33%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -1.80%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +2.19%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +8.75%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.341'
$ws.Cells.Item(13, 5).Value = '  +2.32%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '3.033.62'
$ws.Cells.Item(14, 5).Value = '  +0.00%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '59.096.42'
$ws.Cells.Item(15, 5).Value = '  +1.02%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '21.89'
$ws.Cells.Item(16, 5).Value = '  +6.15%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +3.27%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.581.65'
$ws.Cells.Item(18, 5).Value = '  +0.14%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.61%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '335.15'
$ws.Cells.Item(20, 5).Value = '  +0.20%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.15'
$ws.Cells.Item(21, 5).Value = '  +0.88%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +1.16%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.09%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '64.71'
$ws.Cells.Item(24, 5).Value = '  -2.51%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.444'
$ws.Cells.Item(25, 5).Value = '  +5.04%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.998'
$ws.Cells.Item(26, 5).Value = '  -0.07%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +1.89%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.21'
$ws.Cells.Item(28, 5).Value = '  +1.78%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '0.0₃0778'
$ws.Cells.Item(29, 5).Value = '  +4.82%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.00%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.67'
$ws.Cells.Item(31, 5).Value = '  +2.19%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '160.49'
$ws.Cells.Item(32, 5).Value = '  +4.92%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +1.12%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '18.86'
$ws.Cells.Item(34, 5).Value = '  -0.28%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +2.21%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'SuiNetwork'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.875'
$ws.Cells.Item(36, 5).Value = '  +3.33%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Fetch.AI'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.879'
$ws.Cells.Item(37, 5).Value = '  +7.09%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.13'
$ws.Cells.Item(38, 5).Value = '  +2.88%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '36.74'
$ws.Cells.Item(39, 5).Value = '  -0.95%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +3.92%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '294.68'
$ws.Cells.Item(41, 5).Value = '  +5.61%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +0.61%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +0.10%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.19%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0535'
$ws.Cells.Item(46, 5).Value = '  +1.12%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '10.62'
$ws.Cells.Item(47, 5).Value = '  -0.02%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '125.03'
$ws.Cells.Item(48, 5).Value = '  +14.47%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +2.20%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +1.65%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '18.38'
$ws.Cells.Item(51, 5).Value = '  +2.74%  '
